# Add "2022-Q4" quarterly holding data.
#
# 1) "总计" (summary) sheet: insert a new row 2 for "2022-Q4" (count=9, value=0.06),
#    pushing the existing quarters down by one row.
# 2) Insert a brand-new worksheet named "2022-Q4" right after "总计", containing the
#    per-fund holding detail for that quarter (mirrors the layout used by the other
#    per-quarter sheets).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Update the "总计" summary sheet.
# ---------------------------------------------------------------------------
$totalWs = $wb.Worksheets.Item(1)

# Make room for the new quarter at row 2; existing rows 2.. shift down to 3..
$totalWs.Rows.Item(2).Insert()

# The Insert() above leaves the new A2 with no format at all (it was blank before
# the insert) - clone the number/border/font formatting used by the other index
# cells in column A so it matches (style "bold, centered, thin border").
$totalWs.Range("A3").Copy()
$totalWs.Range("A2").PasteSpecial(-4122)
$totalWs.Application.CutCopyMode = $false

$totalWs.Range("A2").Value = 0
$totalWs.Range("B2").Value = "2022-Q4"
$totalWs.Range("C2").Value = 9
$totalWs.Range("D2").Value = 0.06
$totalWs.Range("B2:D2").Style = "Normal"

$totalWs.Range("A3").Value = 1
$totalWs.Range("B3").Value = "2022-Q3"
$totalWs.Range("C3").Value = 1
$totalWs.Range("D3").Value = 0

$totalWs.Range("A4").Value = 2
$totalWs.Range("B4").Value = "2021-Q4"
$totalWs.Range("C4").Value = 6
$totalWs.Range("D4").Value = 0.26

$totalWs.Range("A5").Value = 3
$totalWs.Range("B5").Value = "2021-Q3"
$totalWs.Range("C5").Value = 5
$totalWs.Range("D5").Value = 0.54

$totalWs.Range("A6").Value = 4
$totalWs.Range("B6").Value = "2021-Q2"
$totalWs.Range("C6").Value = 4
$totalWs.Range("D6").Value = 0.76

# ---------------------------------------------------------------------------
# 2) Insert the new "2022-Q4" detail sheet right after "总计".
# ---------------------------------------------------------------------------
$newWs = $wb.Worksheets.Add($null, $totalWs)
$newWs.Name = "2022-Q4"

$newWs.PageSetup.LeftMargin = 54
$newWs.PageSetup.RightMargin = 54
$newWs.PageSetup.TopMargin = 72
$newWs.PageSetup.BottomMargin = 72
$newWs.PageSetup.HeaderMargin = 36
$newWs.PageSetup.FooterMargin = 36

$newWs.Range("B1").Value = "基金代码"
$newWs.Range("C1").Value = "基金名称"
$newWs.Range("D1").Value = "基金规模"
$newWs.Range("E1").Value = "股票总仓位"
$newWs.Range("F1").Value = "仓位占比"
$newWs.Range("G1").Value = "持有市值(亿元)"
$newWs.Range("H1").Value = "仓位排名"

# Column B (fund code) and D:G (numeric-looking figures) are stored as plain text
# in this workbook's convention, so lead with an apostrophe to stop Excel from
# re-interpreting them (and losing leading zeros) as numbers.
$data = @(
  @(0, "'005460", "银河嘉谊灵活配置混合C",   "'2.65", "'39.19", "'0.98", "'0.0260", 3),
  @(1, "'009619", "博时女性消费主题混合A",   "'0.56", "'72.72", "'3.39", "'0.0190", 7),
  @(2, "'007288", "合煦智远消费主题股票C",   "'0.11", "'83.65", "'5.02", "'0.0055", 3),
  @(3, "'005167", "嘉实润泽量化一年定期开放混合", "'0.56", "'27.25", "'0.58", "'0.0032", 9),
  @(4, "'007287", "合煦智远消费主题股票A",   "'0.03", "'83.65", "'5.02", "'0.0015", 3),
  @(5, "'015921", "申万菱信国证2000指数增强A", "'0.21", "'94.00", "'0.53", "'0.0011", 3),
  @(6, "'009620", "博时女性消费主题混合C",   "'0.03", "'72.72", "'3.39", "'0.0010", 7),
  @(7, "'015922", "申万菱信国证2000指数增强C", "'0.08", "'94.00", "'0.53", "'0.0004", 3),
  @(8, "'005459", "银河嘉谊灵活配置混合A",   "'0.01", "'39.19", "'0.98", "'0.0001", 3)
)

$r = 2
foreach ($row in $data) {
    $newWs.Cells.Item($r, 1).Value = $row[0]
    $newWs.Cells.Item($r, 2).Value = $row[1]
    $newWs.Cells.Item($r, 3).Value = $row[2]
    $newWs.Cells.Item($r, 4).Value = $row[3]
    $newWs.Cells.Item($r, 5).Value = $row[4]
    $newWs.Cells.Item($r, 6).Value = $row[5]
    $newWs.Cells.Item($r, 7).Value = $row[6]
    $newWs.Cells.Item($r, 8).Value = $row[7]
    $r = $r + 1
}

# Re-apply the canonical formatting used on the sibling quarter sheets (bold +
# centered + thin-bordered header row and index column; everything else plain)
# by cloning it from the "2022-Q3" sheet, which sits right after the new sheet.
$fmtSrc = $wb.Worksheets.Item($newWs.Index + 1)
$fmtSrc.Range("B1").Copy()
$newWs.Range("B1:H1").PasteSpecial(-4122)
$fmtSrc.Range("A2").Copy()
$newWs.Range("A2:A10").PasteSpecial(-4122)
$fmtSrc.Range("B2").Copy()
$newWs.Range("B2:G10").PasteSpecial(-4122)
$fmtSrc.Range("H2").Copy()
$newWs.Range("H2:H10").PasteSpecial(-4122)
$newWs.Application.CutCopyMode = $false
